$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("degradation_rates")
$ws.Rows.Item(3).Delete()
